$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the bulk-approval sample rows with fresh placeholder users ---
$ws.Range("A2").Value = "new_user01@yopmail.com"
$ws.Range("B2").Value = 9922222222

$ws.Range("A3").Value = "new_user02@yopmail.com"
$ws.Range("B3").Value = 8520222222

$ws.Range("A4").Value = "new_user03@yopmail.com"
$ws.Range("B4").Value = 8522222333

$ws.Range("A5").Value = "new_user04@yopmail.com"
$ws.Range("B5").Value = 9788555555

$ws.Range("A6").Value = "new_user05@yopmail.com"
$ws.Range("B6").Value = 9555522222

$ws.Range("A7").Value = "new_user06@yopmail.com"
$ws.Range("B7").Value = 9555222222

# --- Rebuild the mailto hyperlinks to match the new email addresses ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:new_user01@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A3:A7"), "mailto:new_user01@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:new_user02@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:new_user03@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:new_user04@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:new_user05@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:new_user06@yopmail.com")

# --- Match the author's final selection state ---
[void]$ws.Range("D4").Select()
